$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename TC6's name: "Yapilacak_is_BirdenFazla_Ekleme" -> "..._ClearCompleted"
$ws.Range("B15").Value = "Yapilacak_is_BirdenFazla_Ekleme_ClearCompleted"

# 2) The old rows 18-20 (TC7's "Todomvc_Sayfası" + "Yapilacak_is_Ekleme" boilerplate
#    steps) get removed -- what used to be row 21 (the Clear-Completed step) slides
#    up to become the new row 18, completing TC6.
$ws.Rows("18:20").Delete()

# 3) Insert 5 fresh rows for the brand new TC7 "Yapilacak_is_Ekleme_All_Active_Completed"
$ws.Rows("19:23").Insert()

# 4) Populate the new TC7 header row
$ws.Range("A19").Value = 7
$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("B19").Value = "Yapilacak_is_Ekleme_All_Active_Completed "
$ws.Range("C19").Value = "* https://todomvc.com/examples/vue/#   adresine girilir"
$ws.Range("D19").Value = "Sayfanın Başarılı şekilde açıldıgı görülür"

# 5) Populate the new TC7 step rows
$ws.Range("C20").Value = "* What needs to be done ? İnput'u içerisine herhangi bir değer girilip Enter'a basılır bu islem 3 defa tekrarlanir"
$ws.Range("D20").Value = "Yazilan 3  değer başarılı şekilde eklendiği görülür,  X item left yazisi görülür , All Active Completed butonları görülür || Eklenen deger active olarak eklendigi gorulur "
$ws.Rows(20).RowHeight = 45

$ws.Range("C21").Value = "* 3 deger icerisinden 1. deger secilir"
$ws.Range("D21").Value = "X items left degerinin azaldigi görülür , yazinin secildigi görülür. Clear Completed butonu aktiflesir"
$ws.Rows(21).RowHeight = 30

$ws.Range("C22").Value = "* Active butonuna tiklanir "
$ws.Range("D22").Value = "Secilmeyen degerlerin kaldigi görülur"

$ws.Range("C23").Value = "* Completed butonuna tiklanir"
$ws.Range("D23").Value = "Secilmeyen degerler ekrandan kaybolur , sadece secilen deger goruntulenir"

# 6) Match the saved view state (scrolled down, C16 selected)
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollRow = 7
